# refactor: changed column names
# The "id" header in column A is renamed to "number".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "number"

# Move the cursor back to the top of the sheet / first data row.
$ws.Range("A2").Select()
